$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing data rows (previously rows 2-9) down by one row to make room
# for the new "Watchdog" row, without a full-row insert (that would also shift
# the unrelated formatted cell at C23 further down, which must stay put).
for ($r = 9; $r -ge 2; $r--) {
    $destRow = $r + 1
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value2
}

# Fill in the newly freed row 2 with the Watchdog message
$ws.Cells.Item(2, 1).Value = "Watchdog"
$ws.Cells.Item(2, 2).Value = 99
$ws.Cells.Item(2, 3).Value = "Watchdog mandato dai client all'host"

# Update the description of ClientDisconnectedMessage, now at row 10
$ws.Cells.Item(10, 3).Value = "Notifica che un client è stato disconnesso"

# Append new row 11 with ClientConnectionLost message
$ws.Cells.Item(11, 1).Value = "ClientConnectionLost"
$ws.Cells.Item(11, 2).Value = 1028
$ws.Cells.Item(11, 3).Value = "Notifica i client che è stata persa la connessione con uno dei client"

# Update selection: select C11 as active cell
$ws.Range("C11").Select()
